$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "weight" values for the date row (D5) and author row (D6)
$ws.Range("D5").Value = 300
$ws.Range("D6").Value = 300

# Move the active selection from B8 to D7
$ws.Range("D7").Select()
